# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 113
$wsExhibit.Range("F9").Value = 1271
$wsExhibit.Range("F11").Value = 1006
$wsExhibit.Range("F12").Value = 10368
$wsExhibit.Range("F16").Value = 1019
$wsExhibit.Range("F18").Value = 11907
$wsExhibit.Range("F19").Value = 12297

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 113
$wsAll.Range("F10").Value = 1271
$wsAll.Range("F12").Value = 1006
$wsAll.Range("F13").Value = 10368
$wsAll.Range("F17").Value = 1019
$wsAll.Range("F19").Value = 11907
$wsAll.Range("F20").Value = 12297
